# Natmi following Dr Hou advice
# Updates the Fbn1-Itgav LR-pair sheet to include the "ECs" sending-cluster
# group (previously only FAPs/M2/sCs were present as sending clusters),
# and refreshes all the computed statistic columns (E:T) for every
# Sending-cluster x Target-cluster combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A:T).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the various expression / specificity metrics.
$data = @(
    @('ECs', 'Fbn1', 'Itgav', 'ECs', 2, 0.6666666666666666, 8.970950666666665, 26.912852, 0.02838798528205506, 0.02838798528205506, 3, 1, 15.35884066666667, 46.076522, 0.1012042817263867, 0.1012042817263867, 137.7834019178604, 1240.050617260744, 0.002872985660129619, 0.002872985660129619),
    @('ECs', 'Fbn1', 'Itgav', 'FAPs', 2, 0.6666666666666666, 8.970950666666665, 26.912852, 0.02838798528205506, 0.02838798528205506, 3, 1, 50.59256466666667, 151.777694, 0.3333704853712116, 0.3333704853712116, 453.8634017248097, 4084.770615523288, 0.009463716432189507, 0.009463716432189507),
    @('ECs', 'Fbn1', 'Itgav', 'M2', 2, 0.6666666666666666, 8.970950666666665, 26.912852, 0.02838798528205506, 0.02838798528205506, 3, 1, 60.37715666666667, 181.13147, 0.397844271305776, 0.397844271305776, 541.640493850271, 4874.76444465244, 0.01129399731837829, 0.01129399731837829),
    @('ECs', 'Fbn1', 'Itgav', 'sCs', 2, 0.6666666666666666, 8.970950666666665, 26.912852, 0.02838798528205506, 0.02838798528205506, 3, 1, 25.43221733333333, 76.296652, 0.1675809615966257, 0.1675809615966258, 228.1511670412782, 2053.360503371503, 0.004757285871357647, 0.004757285871357647),
    @('FAPs', 'Fbn1', 'Itgav', 'ECs', 3, 1, 272.2141316666667, 816.6423950000001, 0.8614037742994388, 0.8614037742994389, 3, 1, 15.35884066666667, 46.076522, 0.1012042817263867, 0.1012042817263867, 4180.893475483354, 37628.04127935019, 0.08717775025437319, 0.08717775025437322),
    @('FAPs', 'Fbn1', 'Itgav', 'FAPs', 3, 1, 272.2141316666667, 816.6423950000001, 0.8614037742994388, 0.8614037742994389, 3, 1, 50.59256466666667, 151.777694, 0.3333704853712116, 0.3333704853712116, 13772.01105952635, 123948.0995357371, 0.2871665943387975, 0.2871665943387975),
    @('FAPs', 'Fbn1', 'Itgav', 'M2', 3, 1, 272.2141316666667, 816.6423950000001, 0.8614037742994388, 0.8614037742994389, 3, 1, 60.37715666666667, 181.13147, 0.397844271305776, 0.397844271305776, 16435.51527451896, 147919.6374706707, 0.3427045568862054, 0.3427045568862054),
    @('FAPs', 'Fbn1', 'Itgav', 'sCs', 3, 1, 272.2141316666667, 816.6423950000001, 0.8614037742994388, 0.8614037742994389, 3, 1, 25.43221733333333, 76.296652, 0.1675809615966257, 0.1675809615966258, 6923.008957751282, 62307.08061976154, 0.1443548728200627, 0.1443548728200628),
    @('M2', 'Fbn1', 'Itgav', 'ECs', 3, 1, 0.405826, 1.217478, 0.001284209772536402, 0.001284209772536402, 3, 1, 15.35884066666667, 46.076522, 0.1012042817263867, 0.1012042817263867, 6.233016872390666, 56.097151851516, 0.0001299675276155529, 0.0001299675276155529),
    @('M2', 'Fbn1', 'Itgav', 'FAPs', 3, 1, 0.405826, 1.217478, 0.001284209772536402, 0.001284209772536402, 3, 1, 50.59256466666667, 151.777694, 0.3333704853712116, 0.3333704853712116, 20.53177814841467, 184.786003335732, 0.0004281176351889135, 0.0004281176351889135),
    @('M2', 'Fbn1', 'Itgav', 'M2', 3, 1, 0.405826, 1.217478, 0.001284209772536402, 0.001284209772536402, 3, 1, 60.37715666666667, 181.13147, 0.397844271305776, 0.397844271305776, 24.50261998140667, 220.52357983266, 0.0005109155011585011, 0.0005109155011585011),
    @('M2', 'Fbn1', 'Itgav', 'sCs', 3, 1, 0.405826, 1.217478, 0.001284209772536402, 0.001284209772536402, 3, 1, 25.43221733333333, 76.296652, 0.1675809615966257, 0.1675809615966258, 10.32105503151733, 92.889495283656, 0.0002152091085734342, 0.0002152091085734343),
    @('sCs', 'Fbn1', 'Itgav', 'ECs', 3, 1, 34.42132633333333, 103.263979, 0.1089240306459696, 0.1089240306459696, 3, 1, 15.35884066666667, 46.076522, 0.1012042817263867, 0.1012042817263867, 528.6716666890042, 4758.045000201038, 0.01102357828426828, 0.01102357828426828),
    @('sCs', 'Fbn1', 'Itgav', 'FAPs', 3, 1, 34.42132633333333, 103.263979, 0.1089240306459696, 0.1089240306459696, 3, 1, 50.59256466666667, 151.777694, 0.3333704853712116, 0.3333704853712116, 1741.463178431603, 15673.16860588442, 0.0363120569650356, 0.0363120569650356),
    @('sCs', 'Fbn1', 'Itgav', 'M2', 3, 1, 34.42132633333333, 103.263979, 0.1089240306459696, 0.1089240306459696, 3, 1, 60.37715666666667, 181.13147, 0.397844271305776, 0.397844271305776, 2078.261812702126, 18704.35631431913, 0.04333480160003378, 0.04333480160003378),
    @('sCs', 'Fbn1', 'Itgav', 'sCs', 3, 1, 34.42132633333333, 103.263979, 0.1089240306459696, 0.1089240306459696, 3, 1, 25.43221733333333, 76.296652, 0.1675809615966257, 0.1675809615966258, 875.410652210923, 7878.695869898307, 0.01825359379663192, 0.01825359379663192)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
